$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cases" query (row 2 / cell B2) dropped the trailing Cohort column
# from its RETURN clause -- update the cell text accordingly. (Samples /
# Files tab queries, B3 & B4, are unchanged in content.)
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (co:cohort)<-[*]-(c)
WHERE diag.disease_term IN ['T Cell Lymphoma']
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value2 = $casesQuery

# Selection moved from B4 to B2, and the view scrolled back so row 1 is
# visible at the top again (topLeftCell reset to A1).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$null = $ws.Range("B2").Select()
